{"js": "// Replace each two-digit multiplication expression's old text with its new value.\n// Each \"old\" string is unique across the document, so an exact, case-sensitive,\n// non-wildcard search finds exactly one match per entry.\nconst replacements = [\n  [\"63\u00d748=3024\", \"25\u00d711=275\"],\n  [\"17\u00d732=544\", \"55\u00d724=1320\"],\n  [\"24\u00d733=792\", \"99\u00d711=1089\"],\n  [\"47\u00d748=2256\", \"59\u00d787=5133\"],\n  [\"59\u00d721=1239\", \"78\u00d773=5694\"],\n  [\"70\u00d746=3220\", \"39\u00d778=3042\"],\n  [\"68\u00d798=6664\", \"17\u00d749=833\"],\n  [\"15\u00d759=885\", \"15\u00d758=870\"],\n  [\"99\u00d760=5940\", \"37\u00d774=2738\"],\n  [\"67\u00d720=1340\", \"27\u00d743=1161\"],\n  [\"83\u00d732=2656\", \"35\u00d746=1610\"],\n  [\"25\u00d762=1550\", \"50\u00d741=2050\"],\n  [\"93\u00d740=3720\", \"22\u00d755=1210\"],\n  [\"28\u00d759=1652\", \"66\u00d766=4356\"],\n  [\"46\u00d752=2392\", \"59\u00d734=2006\"],\n  [\"74\u00d736=2664\", \"44\u00d713=572\"],\n  [\"87\u00d714=1218\", \"46\u00d781=3726\"],\n  [\"11\u00d755=605\", \"43\u00d745=1935\"],\n  [\"96\u00d758=5568\", \"24\u00d763=1512\"],\n  [\"46\u00d756=2576\", \"68\u00d754=3672\"],\n  [\"60\u00d725=1500\", \"66\u00d745=2970\"],\n  [\"40\u00d719=760\", \"54\u00d753=2862\"],\n  [\"23\u00d755=1265\", \"39\u00d741=1599\"],\n  [\"73\u00d722=1606\", \"46\u00d773=3358\"],\n  [\"36\u00d776=2736\", \"62\u00d715=930\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression's old text with its new value.\n# Each \"old\" string is unique across the document, so Find/Replace with exact,\n# case-sensitive, non-wildcard matching touches exactly one run per entry.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '63\u00d748=3024'; New = '25\u00d711=275' },\n    @{ Old = '17\u00d732=544'; New = '55\u00d724=1320' },\n    @{ Old = '24\u00d733=792'; New = '99\u00d711=1089' },\n    @{ Old = '47\u00d748=2256'; New = '59\u00d787=5133' },\n    @{ Old = '59\u00d721=1239'; New = '78\u00d773=5694' },\n    @{ Old = '70\u00d746=3220'; New = '39\u00d778=3042' },\n    @{ Old = '68\u00d798=6664'; New = '17\u00d749=833' },\n    @{ Old = '15\u00d759=885'; New = '15\u00d758=870' },\n    @{ Old = '99\u00d760=5940'; New = '37\u00d774=2738' },\n    @{ Old = '67\u00d720=1340'; New = '27\u00d743=1161' },\n    @{ Old = '83\u00d732=2656'; New = '35\u00d746=1610' },\n    @{ Old = '25\u00d762=1550'; New = '50\u00d741=2050' },\n    @{ Old = '93\u00d740=3720'; New = '22\u00d755=1210' },\n    @{ Old = '28\u00d759=1652'; New = '66\u00d766=4356' },\n    @{ Old = '46\u00d752=2392'; New = '59\u00d734=2006' },\n    @{ Old = '74\u00d736=2664'; New = '44\u00d713=572' },\n    @{ Old = '87\u00d714=1218'; New = '46\u00d781=3726' },\n    @{ Old = '11\u00d755=605'; New = '43\u00d745=1935' },\n    @{ Old = '96\u00d758=5568'; New = '24\u00d763=1512' },\n    @{ Old = '46\u00d756=2576'; New = '68\u00d754=3672' },\n    @{ Old = '60\u00d725=1500'; New = '66\u00d745=2970' },\n    @{ Old = '40\u00d719=760'; New = '54\u00d753=2862' },\n    @{ Old = '23\u00d755=1265'; New = '39\u00d741=1599' },\n    @{ Old = '73\u00d722=1606'; New = '46\u00d773=3358' },\n    @{ Old = '36\u00d776=2736'; New = '62\u00d715=930' }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $found = $find.Execute(\n        $r.Old,      # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $r.New,      # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"No match found for '$($r.Old)'\"\n    }\n}\n"}
